$d = $word.ActiveDocument

# 1. Fix "service mysql start" -> "service mysqld start": add a new run "d" after "mysql"
$d.Content.Find.Execute("service mysql start", $true, $false, $false, $false, $false, $true, 1, $false, "service mysqld start", 2) | Out-Null

# 2. Append the new "Allowing to write files from filezilla" section before the
#    trailing empty paragraph, using the literal OOXML captured from the source edit.
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$insertionRange = $lastPara.Range
$insertionRange.Collapse(1)

$newContentXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>    <w:p/>
    <w:p>
      <w:r>
        <w:lastRenderedPageBreak/>
        <w:t>Allowing to write files from filezilla</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
        <w:textAlignment w:val="baseline"/>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
          <w:color w:val="242729"/>
          <w:sz w:val="23"/>
          <w:szCs w:val="23"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
          <w:color w:val="242729"/>
          <w:sz w:val="23"/>
          <w:szCs w:val="23"/>
        </w:rPr>
        <w:t>To allow user </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Courier New"/>
          <w:color w:val="242729"/>
          <w:sz w:val="20"/>
        </w:rPr>
        <w:t>ec2-user</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
          <w:color w:val="242729"/>
          <w:sz w:val="23"/>
          <w:szCs w:val="23"/>
        </w:rPr>
        <w:t> (Amazon AWS) write access to the public web directory (/</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
          <w:color w:val="242729"/>
          <w:sz w:val="23"/>
          <w:szCs w:val="23"/>
        </w:rPr>
        <w:t>var</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
          <w:color w:val="242729"/>
          <w:sz w:val="23"/>
          <w:szCs w:val="23"/>
        </w:rPr>
        <w:t>/www/html)</w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
          <w:color w:val="242729"/>
          <w:sz w:val="23"/>
          <w:szCs w:val="23"/>
        </w:rPr>
        <w:t>,</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
          <w:color w:val="242729"/>
          <w:sz w:val="23"/>
          <w:szCs w:val="23"/>
        </w:rPr>
        <w:br/>
        <w:t>enter this command via Putty or Terminal, as the root user </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Courier New"/>
          <w:color w:val="242729"/>
          <w:sz w:val="20"/>
        </w:rPr>
        <w:t>sudo</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
          <w:color w:val="242729"/>
          <w:sz w:val="23"/>
          <w:szCs w:val="23"/>
        </w:rPr>
        <w:t>:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:shd w:val="clear" w:color="auto" w:fill="EFF0F1"/>
        <w:tabs>
          <w:tab w:val="left" w:pos="916"/>
          <w:tab w:val="left" w:pos="1832"/>
          <w:tab w:val="left" w:pos="2748"/>
          <w:tab w:val="left" w:pos="3664"/>
          <w:tab w:val="left" w:pos="4580"/>
          <w:tab w:val="left" w:pos="5496"/>
          <w:tab w:val="left" w:pos="6412"/>
          <w:tab w:val="left" w:pos="7328"/>
          <w:tab w:val="left" w:pos="8244"/>
          <w:tab w:val="left" w:pos="9160"/>
          <w:tab w:val="left" w:pos="10076"/>
          <w:tab w:val="left" w:pos="10992"/>
          <w:tab w:val="left" w:pos="11908"/>
          <w:tab w:val="left" w:pos="12824"/>
          <w:tab w:val="left" w:pos="13740"/>
          <w:tab w:val="left" w:pos="14656"/>
        </w:tabs>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
        <w:textAlignment w:val="baseline"/>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Courier New"/>
          <w:color w:val="242729"/>
          <w:sz w:val="20"/>
        </w:rPr>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Courier New"/>
          <w:color w:val="242729"/>
          <w:sz w:val="20"/>
        </w:rPr>
        <w:t>chown</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Courier New"/>
          <w:color w:val="242729"/>
          <w:sz w:val="20"/>
        </w:rPr>
        <w:t xml:space="preserve"> -R ec2-user /</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Courier New"/>
          <w:color w:val="242729"/>
          <w:sz w:val="20"/>
        </w:rPr>
        <w:t>var</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Courier New"/>
          <w:color w:val="242729"/>
          <w:sz w:val="20"/>
        </w:rPr>
        <w:t>/www/html</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        <w:spacing w:after="240" w:line="240" w:lineRule="auto"/>
        <w:textAlignment w:val="baseline"/>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
          <w:color w:val="242729"/>
          <w:sz w:val="23"/>
          <w:szCs w:val="23"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
          <w:color w:val="242729"/>
          <w:sz w:val="23"/>
          <w:szCs w:val="23"/>
        </w:rPr>
        <w:t>Make sure permissions on that entire folder were correct:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:shd w:val="clear" w:color="auto" w:fill="EFF0F1"/>
        <w:tabs>
          <w:tab w:val="left" w:pos="916"/>
          <w:tab w:val="left" w:pos="1832"/>
          <w:tab w:val="left" w:pos="2748"/>
          <w:tab w:val="left" w:pos="3664"/>
          <w:tab w:val="left" w:pos="4580"/>
          <w:tab w:val="left" w:pos="5496"/>
          <w:tab w:val="left" w:pos="6412"/>
          <w:tab w:val="left" w:pos="7328"/>
          <w:tab w:val="left" w:pos="8244"/>
          <w:tab w:val="left" w:pos="9160"/>
          <w:tab w:val="left" w:pos="10076"/>
          <w:tab w:val="left" w:pos="10992"/>
          <w:tab w:val="left" w:pos="11908"/>
          <w:tab w:val="left" w:pos="12824"/>
          <w:tab w:val="left" w:pos="13740"/>
          <w:tab w:val="left" w:pos="14656"/>
        </w:tabs>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
        <w:textAlignment w:val="baseline"/>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Courier New"/>
          <w:color w:val="242729"/>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
        </w:rPr>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Courier New"/>
          <w:color w:val="242729"/>
          <w:sz w:val="20"/>
        </w:rPr>
        <w:t>chmod</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Courier New"/>
          <w:color w:val="242729"/>
          <w:sz w:val="20"/>
        </w:rPr>
        <w:t xml:space="preserve"> -R 755 /</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Courier New"/>
          <w:color w:val="242729"/>
          <w:sz w:val="20"/>
        </w:rPr>
        <w:t>var</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Courier New"/>
          <w:color w:val="242729"/>
          <w:sz w:val="20"/>
        </w:rPr>
        <w:t>/www/html</w:t>
      </w:r>
    </w:p>
    <w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertionRange.InsertXML($newContentXml)

Write-Output "applied edits"
